$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Cells.Item(38, 8).Value = 542.2857
$ws.Cells.Item(38, 9).Value = 698.5
$ws.Cells.Item(38, 10).Value = 479.8
$ws.Cells.Item(38, 11).Value = 2095.5
$ws.Cells.Item(38, 12).Value = 1439.4
$ws.Cells.Item(38, 13).Value = -1723.5
$ws.Cells.Item(38, 14).Value = -2183.4

# Row 40
$ws.Cells.Item(40, 8).Value = 2045
$ws.Cells.Item(40, 9).Value = 1984.45
$ws.Cells.Item(40, 10).Value = 2138.1538
$ws.Cells.Item(40, 11).Value = 1984.45
$ws.Cells.Item(40, 12).Value = 2138.1538
$ws.Cells.Item(40, 13).Value = -1809.45
$ws.Cells.Item(40, 14).Value = -2488.1538

# Row 43
$ws.Cells.Item(43, 8).Value = 559.6667
$ws.Cells.Item(43, 9).Value = 495.33334
$ws.Cells.Item(43, 10).Value = 624
$ws.Cells.Item(43, 11).Value = 495.33334
$ws.Cells.Item(43, 12).Value = 624
$ws.Cells.Item(43, 13).Value = -426.33334
$ws.Cells.Item(43, 14).Value = -762

# Row 44
$ws.Cells.Item(44, 8).Value = 49000
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 49000
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 49000
$ws.Cells.Item(44, 14).Value = -49924

# Row 113
$ws.Cells.Item(113, 8).Value = 111080.5
$ws.Cells.Item(113, 9).Value = 178800.83
$ws.Cells.Item(113, 10).Value = 9500
$ws.Cells.Item(113, 11).Value = 178800.83
$ws.Cells.Item(113, 12).Value = 9500
$ws.Cells.Item(113, 13).Value = -175546.83
$ws.Cells.Item(113, 14).Value = -16008

# Row 132
$ws.Cells.Item(132, 8).Value = 326554.12
$ws.Cells.Item(132, 9).Value = 419876.03
$ws.Cells.Item(132, 10).Value = 55920.6
$ws.Cells.Item(132, 11).Value = 1259628.09
$ws.Cells.Item(132, 12).Value = 167761.8
$ws.Cells.Item(132, 13).Value = -1257098.09
$ws.Cells.Item(132, 14).Value = -172821.8

# Row 133
$ws.Cells.Item(133, 8).Value = 16097
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 16097
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 16097
$ws.Cells.Item(133, 14).Value = -26217

$ws = $wb.Worksheets.Item("ARM")
# Row 31
$ws.Cells.Item(31, 8).Value = 9109
$ws.Cells.Item(31, 9).Value = 3136.25
$ws.Cells.Item(31, 10).Value = 33000
$ws.Cells.Item(31, 11).Value = 3136.25
$ws.Cells.Item(31, 12).Value = 33000
$ws.Cells.Item(31, 13).Value = -2842.25
$ws.Cells.Item(31, 14).Value = -33588

# Row 32
$ws.Cells.Item(32, 8).Value = 14573.65
$ws.Cells.Item(32, 9).Value = 1505.04
$ws.Cells.Item(32, 10).Value = 210602.8
$ws.Cells.Item(32, 11).Value = 1505.04
$ws.Cells.Item(32, 12).Value = 210602.8
$ws.Cells.Item(32, 13).Value = -1218.04
$ws.Cells.Item(32, 14).Value = -211176.8

# Row 122
$ws.Cells.Item(122, 8).Value = 3007
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 3007
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).ClearContents()
$ws.Cells.Item(122, 13).Value = 9021
$ws.Cells.Item(122, 14).Value = -13921

# Row 132
$ws.Cells.Item(132, 8).Value = 2883.4443
$ws.Cells.Item(132, 9).Value = 2492.5
$ws.Cells.Item(132, 10).Value = 4838.1665
$ws.Cells.Item(132, 11).Value = 7477.5
$ws.Cells.Item(132, 12).Value = 14514.4995
$ws.Cells.Item(132, 13).Value = -4947.5
$ws.Cells.Item(132, 14).Value = -19574.4995

# Row 133
$ws.Cells.Item(133, 8).Value = 38844.445
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 38844.445
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 38844.445
$ws.Cells.Item(133, 14).Value = -43904.445

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1530.1765
$ws.Cells.Item(20, 9).Value = 1476.0834
$ws.Cells.Item(20, 10).Value = 1660
$ws.Cells.Item(20, 11).Value = 1476.0834
$ws.Cells.Item(20, 12).Value = 1660
$ws.Cells.Item(20, 13).Value = -1229.0834
$ws.Cells.Item(20, 14).Value = -2154

# Row 59
$ws.Cells.Item(59, 8).Value = 47166.668
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 10).Value = 47166.668
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 12).Value = 47166.668
$ws.Cells.Item(59, 14).Value = -48860.668

# Row 122
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).ClearContents()
$ws.Cells.Item(122, 14).Value = 0

# Row 124
$ws.Cells.Item(124, 8).Value = 44160
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 44160
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 44160
$ws.Cells.Item(124, 14).Value = -53980

# Row 125
$ws.Cells.Item(125, 8).Value = 40000
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 40000
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 40000
$ws.Cells.Item(125, 14).Value = -49840

# Row 126
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).ClearContents()

# Row 133
$ws.Cells.Item(133, 8).Value = 37775.75
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 37775.75
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 37775.75
$ws.Cells.Item(133, 14).Value = -47895.75

# Row 134
$ws.Cells.Item(134, 8).Value = 3347.44
$ws.Cells.Item(134, 9).Value = 2659.238
$ws.Cells.Item(134, 10).Value = 6960.5
$ws.Cells.Item(134, 11).Value = 7977.714
$ws.Cells.Item(134, 12).Value = 20881.5
$ws.Cells.Item(134, 13).Value = -5442.714
$ws.Cells.Item(134, 14).Value = -25951.5

$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Cells.Item(20, 8).Value = 49749.5
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 49749.5
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 49749.5
$ws.Cells.Item(20, 14).Value = -50221.5

# Row 30
$ws.Cells.Item(30, 8).Value = 49749.5
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 49749.5
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 49749.5
$ws.Cells.Item(30, 14).Value = -49931.5

# Row 31
$ws.Cells.Item(31, 8).Value = 1019.0833
$ws.Cells.Item(31, 9).Value = 887.85297
$ws.Cells.Item(31, 10).Value = 3250
$ws.Cells.Item(31, 11).Value = 887.85297
$ws.Cells.Item(31, 12).Value = 3250
$ws.Cells.Item(31, 13).Value = -592.85297
$ws.Cells.Item(31, 14).Value = -3840

# Row 33
$ws.Cells.Item(33, 8).Value = 9385.143
$ws.Cells.Item(33, 9).Value = 9385.143
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 11).Value = 9385.143
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 13).Value = -9006.143

# Row 34
$ws.Cells.Item(34, 8).Value = 1019.0833
$ws.Cells.Item(34, 9).Value = 887.85297
$ws.Cells.Item(34, 10).Value = 3250
$ws.Cells.Item(34, 11).Value = 887.85297
$ws.Cells.Item(34, 12).Value = 3250
$ws.Cells.Item(34, 13).Value = -685.85297
$ws.Cells.Item(34, 14).Value = -3654

# Row 86
$ws.Cells.Item(86, 8).Value = 41668616
$ws.Cells.Item(86, 9).Value = 71430376
$ws.Cells.Item(86, 10).Value = 2158
$ws.Cells.Item(86, 11).Value = 71430376
$ws.Cells.Item(86, 12).Value = 2158
$ws.Cells.Item(86, 13).Value = -71429253
$ws.Cells.Item(86, 14).Value = -4404

# Row 89
$ws.Cells.Item(89, 8).Value = 41668616
$ws.Cells.Item(89, 9).Value = 71430376
$ws.Cells.Item(89, 10).Value = 2158
$ws.Cells.Item(89, 11).Value = 357151880
$ws.Cells.Item(89, 12).Value = 10790
$ws.Cells.Item(89, 13).Value = -357146264
$ws.Cells.Item(89, 14).Value = -22022

# Row 94
$ws.Cells.Item(94, 8).Value = 1269.0834
$ws.Cells.Item(94, 9).Value = 702.5454999999999
$ws.Cells.Item(94, 10).Value = 1748.4615
$ws.Cells.Item(94, 11).Value = 702.5454999999999
$ws.Cells.Item(94, 12).Value = 1748.4615
$ws.Cells.Item(94, 13).Value = -251.5454999999999
$ws.Cells.Item(94, 14).Value = -2650.4615

# Row 128
$ws.Cells.Item(128, 8).Value = 49749.5
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 49749.5
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 49749.5
$ws.Cells.Item(128, 14).Value = -59709.5

# Row 130
$ws.Cells.Item(130, 8).Value = 48500
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 48500
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 48500
$ws.Cells.Item(130, 14).Value = -58540

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Cells.Item(4, 8).Value = 31089.908
$ws.Cells.Item(4, 9).Value = 97
$ws.Cells.Item(4, 10).Value = 56917.332
$ws.Cells.Item(4, 11).Value = 291
$ws.Cells.Item(4, 12).Value = 170751.996
$ws.Cells.Item(4, 13).Value = -179
$ws.Cells.Item(4, 14).Value = -170975.996

# Row 129
$ws.Cells.Item(129, 8).Value = 1255.7368
$ws.Cells.Item(129, 9).Value = 439.9
$ws.Cells.Item(129, 10).Value = 2162.2222
$ws.Cells.Item(129, 11).Value = 1319.7
$ws.Cells.Item(129, 12).Value = 6486.6666
$ws.Cells.Item(129, 13).Value = 3680.3
$ws.Cells.Item(129, 14).Value = -16486.6666

# Row 131
$ws.Cells.Item(131, 8).Value = 1438.2142
$ws.Cells.Item(131, 9).Value = 412
$ws.Cells.Item(131, 10).Value = 1609.25
$ws.Cells.Item(131, 11).Value = 1236
$ws.Cells.Item(131, 12).Value = 4827.75
$ws.Cells.Item(131, 13).Value = 3804
$ws.Cells.Item(131, 14).Value = -14907.75

# Row 133
$ws.Cells.Item(133, 8).Value = 9500
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 9500
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).ClearContents()
$ws.Cells.Item(133, 13).Value = 28500
$ws.Cells.Item(133, 14).Value = -38620

# Row 134
$ws.Cells.Item(134, 8).Value = 5185.727
$ws.Cells.Item(134, 9).Value = 3130
$ws.Cells.Item(134, 10).Value = 8783.25
$ws.Cells.Item(134, 11).Value = 9390
$ws.Cells.Item(134, 12).Value = 26349.75
$ws.Cells.Item(134, 13).Value = -4320
$ws.Cells.Item(134, 14).Value = -36489.75

# Row 136
$ws.Cells.Item(136, 8).Value = 2739.697
$ws.Cells.Item(136, 9).Value = 2443.3333
$ws.Cells.Item(136, 10).Value = 2769.3333
$ws.Cells.Item(136, 11).Value = 7329.999899999999
$ws.Cells.Item(136, 12).Value = 8307.999899999999
$ws.Cells.Item(136, 13).Value = -2229.999899999999
$ws.Cells.Item(136, 14).Value = -18507.9999

# Row 137
$ws.Cells.Item(137, 8).Value = 6737550
$ws.Cells.Item(137, 9).Value = 20002134
$ws.Cells.Item(137, 10).Value = 105257.6
$ws.Cells.Item(137, 11).Value = 60006402
$ws.Cells.Item(137, 12).Value = 315772.8
$ws.Cells.Item(137, 13).Value = -60001302
$ws.Cells.Item(137, 14).Value = -325972.8

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 1588926.8
$ws.Cells.Item(122, 9).Value = 3704936
$ws.Cells.Item(122, 10).Value = 1919.75
$ws.Cells.Item(122, 11).Value = 11114808
$ws.Cells.Item(122, 12).Value = 5759.25
$ws.Cells.Item(122, 13).Value = -11112358
$ws.Cells.Item(122, 14).Value = -10659.25

# Row 124
$ws.Cells.Item(124, 8).Value = 60999.668
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 60999.668
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 60999.668
$ws.Cells.Item(124, 14).Value = -70819.66800000001

# Row 132
$ws.Cells.Item(132, 8).Value = 2856.9038
$ws.Cells.Item(132, 9).Value = 2412.276
$ws.Cells.Item(132, 10).Value = 3417.5217
$ws.Cells.Item(132, 11).Value = 7236.828
$ws.Cells.Item(132, 12).Value = 10252.5651
$ws.Cells.Item(132, 13).Value = -4706.828
$ws.Cells.Item(132, 14).Value = -15312.5651

# Row 138
$ws.Cells.Item(138, 8).Value = 49208.89
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 49208.89
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 49208.89
$ws.Cells.Item(138, 14).Value = -59488.89

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 8688.308000000001
$ws.Cells.Item(22, 9).Value = 467.33334
$ws.Cells.Item(22, 10).Value = 11154.6
$ws.Cells.Item(22, 11).Value = 467.33334
$ws.Cells.Item(22, 12).Value = 11154.6
$ws.Cells.Item(22, 13).Value = -172.33334
$ws.Cells.Item(22, 14).Value = -11744.6

# Row 27
$ws.Cells.Item(27, 8).Value = 8688.308000000001
$ws.Cells.Item(27, 9).Value = 467.33334
$ws.Cells.Item(27, 10).Value = 11154.6
$ws.Cells.Item(27, 11).Value = 467.33334
$ws.Cells.Item(27, 12).Value = 11154.6
$ws.Cells.Item(27, 13).Value = -360.33334
$ws.Cells.Item(27, 14).Value = -11368.6

# Row 46
$ws.Cells.Item(46, 8).Value = 1118.8125
$ws.Cells.Item(46, 9).Value = 1040
$ws.Cells.Item(46, 10).Value = 1250.1666
$ws.Cells.Item(46, 11).Value = 1040
$ws.Cells.Item(46, 12).Value = 1250.1666
$ws.Cells.Item(46, 13).Value = -852
$ws.Cells.Item(46, 14).Value = -1626.1666

# Row 99
$ws.Cells.Item(99, 8).Value = 24000
$ws.Cells.Item(99, 9).Value = 16000
$ws.Cells.Item(99, 10).Value = 32000
$ws.Cells.Item(99, 11).Value = 16000
$ws.Cells.Item(99, 12).Value = 32000
$ws.Cells.Item(99, 13).Value = -13005
$ws.Cells.Item(99, 14).Value = -37990

# Row 132
$ws.Cells.Item(132, 8).Value = 4240.9575
$ws.Cells.Item(132, 9).Value = 3666.3142
$ws.Cells.Item(132, 10).Value = 5917
$ws.Cells.Item(132, 11).Value = 10998.9426
$ws.Cells.Item(132, 12).Value = 17751
$ws.Cells.Item(132, 13).Value = -8468.942599999998
$ws.Cells.Item(132, 14).Value = -22811

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Cells.Item(126, 8).Value = 85266.75
$ws.Cells.Item(126, 9).Value = 92872.82000000001
$ws.Cells.Item(126, 10).Value = 1600
$ws.Cells.Item(126, 11).Value = 278618.46
$ws.Cells.Item(126, 12).Value = 4800
$ws.Cells.Item(126, 13).Value = -276148.46
$ws.Cells.Item(126, 14).Value = -9740

# Row 136
$ws.Cells.Item(136, 8).Value = 12384466
$ws.Cells.Item(136, 9).Value = 17597846
$ws.Cells.Item(136, 10).Value = 2689
$ws.Cells.Item(136, 11).Value = 52793538
$ws.Cells.Item(136, 12).Value = 8067
$ws.Cells.Item(136, 13).Value = -52790988
$ws.Cells.Item(136, 14).Value = -13167
